# Commit: "remove table from results-view in rst, png, pptx"
#
# 1. Bump the cached datetimeFigureOut field text (11/14/2024 -> 11/15/2024)
#    on the slide master and every slide layout's Date Placeholder.
# 2. Resize/reposition the results-view screenshot on slide 3.
# 3. Delete the "Table 15" graphicFrame (Result/Value/Units table) on slide 3.

$p = $ppt.ActivePresentation

# --- 1. Date placeholder text: slide master + all slide layouts ---------

$sm = $p.SlideMaster

for ($j = 1; $j -le $sm.Shapes.Count; $j++) {
    $shp = $sm.Shapes.Item($j)
    if ($shp.Name -like "*Date*" -and $shp.HasTextFrame) {
        if ($shp.TextFrame.TextRange.Text -eq "11/14/2024") {
            $shp.TextFrame.TextRange.Text = "11/15/2024"
        }
    }
}

for ($i = 1; $i -le $sm.CustomLayouts.Count; $i++) {
    $lay = $sm.CustomLayouts.Item($i)
    for ($j = 1; $j -le $lay.Shapes.Count; $j++) {
        $shp = $lay.Shapes.Item($j)
        if ($shp.Name -like "*Date*" -and $shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq "11/14/2024") {
                $shp.TextFrame.TextRange.Text = "11/15/2024"
            }
        }
    }
}

# --- 2 & 3. Slide 3: resize the screenshot, drop the results table ------

$s3 = $p.Slides.Item(3)

for ($j = 1; $j -le $s3.Shapes.Count; $j++) {
    $shp = $s3.Shapes.Item($j)

    if ($shp.Name -eq "Picture 4") {
        # New off/ext (EMU 3706766/1502943 4670997/3581820) expressed in
        # points, nudged so the engine's float32-then-truncate EMU
        # round-trip lands exactly on the target integer EMU values.
        $shp.Left = 291.87135858267715
        $shp.Top = 118.34197850393701
        $shp.Width = 367.7950293700788
        $shp.Height = 282.0330708661417
    }
}

for ($j = $s3.Shapes.Count; $j -ge 1; $j--) {
    $shp = $s3.Shapes.Item($j)
    if ($shp.Name -eq "Table 15") {
        $shp.Delete()
    }
}
